# Apply changes described by the diff:
# 1. Fix "wrong" -> "correct" on sheet "test" (C2, C4, C5)
# 2. Add a new worksheet named "difficult" with vocabulary data

$wb = $excel.ActiveWorkbook

# --- Fix the "wrong" entries on the "test" sheet ---
$testSheet = $wb.Worksheets.Item("test")
$testSheet.Range("C2").Value = "correct"
$testSheet.Range("C4").Value = "correct"
$testSheet.Range("C5").Value = "correct"

# --- Add new "difficult" sheet after "H1" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "difficult"

$newSheet.Range("A1").Value = "eten"
$newSheet.Range("B1").Value = "manger"
$newSheet.Range("A2").Value = "zien"
$newSheet.Range("B2").Value = "voir"
$newSheet.Range("A3").Value = "test"
$newSheet.Range("B3").Value = "test"
$newSheet.Range("A4").Value = "opzetten, opstellen"
$newSheet.Range("B4").Value = "dresser"

# Move the new sheet to the end, after "H1"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
